$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the student label in A15: "14 Nombre Apellido1 Apellido2" -> "15 Nombre Apellido1 Apellido2"
$ws.Range("A15").Value = "15 Nombre Apellido1 Apellido2"

# Update the grades for that row
$ws.Range("B15").Value = 10
$ws.Range("C15").Value = 9.5

# Update selection / view state to match the saved workbook
$ws.Range("C16").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
